# Auto-generated Excel COM-interop script
# Applies per-cell value updates (and clears) to match the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 2058.6
$ws.Range("I6").Value = 30
$ws.Range("K6").Value = 90
$ws.Range("M6").Value = 22
# Row 15
$ws.Range("H15").Value = 1145.6316
$ws.Range("I15").Value = 1145.6316
$ws.Range("K15").Value = 3436.8948
$ws.Range("M15").Value = -3267.8948
# Row 16
$ws.Range("H16").Value = 5000
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5460
# Row 18
$ws.Range("H18").Value = 9245.477000000001
$ws.Range("I18").Value = 2497.25
$ws.Range("J18").Value = 13398.23
$ws.Range("K18").Value = 2497.25
$ws.Range("L18").Value = 13398.23
$ws.Range("M18").Value = -2213.25
$ws.Range("N18").Value = -13966.23
# Row 21
$ws.Range("H21").Value = 18152.125
$ws.Range("I21").Value = 18152.125
$ws.Range("K21").Value = 18152.125
$ws.Range("M21").Value = -17684.125
# Row 23
$ws.Range("H23").Value = 18152.125
$ws.Range("I23").Value = 18152.125
$ws.Range("K23").Value = 18152.125
$ws.Range("M23").Value = -17918.125
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 41
$ws.Range("H41").Value = 601
$ws.Range("J41").Value = 730
$ws.Range("L41").Value = 730
$ws.Range("N41").Value = -1610
# Row 96
$ws.Range("H96").Value = 1337.9286
$ws.Range("I96").Value = 1285.2858
$ws.Range("J96").Value = 1390.5714
$ws.Range("K96").Value = 3855.8574
$ws.Range("L96").Value = 4171.7142
$ws.Range("M96").Value = -2482.8574
$ws.Range("N96").Value = -6917.7142
# Row 132
$ws.Range("H132").Value = 1214.9231
$ws.Range("I132").Value = 1241.2084
$ws.Range("K132").Value = 3723.6252
$ws.Range("M132").Value = -1193.6252

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2326339
$ws.Range("I2").Value = 2326339
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2326339
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2326226
$ws.Range("N2").ClearContents()
# Row 32
$ws.Range("H32").Value = 3685.93
$ws.Range("I32").Value = 1984.2826
$ws.Range("J32").Value = 10801.909
$ws.Range("K32").Value = 1984.2826
$ws.Range("L32").Value = 10801.909
$ws.Range("M32").Value = -1697.2826
$ws.Range("N32").Value = -11375.909
# Row 61
$ws.Range("H61").Value = 5937
$ws.Range("I61").Value = 3200.5557
$ws.Range("J61").Value = 9455.286
$ws.Range("K61").Value = 3200.5557
$ws.Range("L61").Value = 9455.286
$ws.Range("M61").Value = -2988.5557
$ws.Range("N61").Value = -9879.286
# Row 63
$ws.Range("H63").Value = 7296
$ws.Range("I63").Value = 7296
$ws.Range("K63").Value = 7296
$ws.Range("M63").Value = -6610
# Row 66
$ws.Range("H66").Value = 7296
$ws.Range("I66").Value = 7296
$ws.Range("K66").Value = 36480
$ws.Range("M66").Value = -33048
# Row 82
$ws.Range("H82").Value = 38333
$ws.Range("J82").Value = 38333
$ws.Range("L82").Value = 38333
$ws.Range("N82").Value = -39055
# Row 85
$ws.Range("H85").Value = 38333
$ws.Range("J85").Value = 38333
$ws.Range("L85").Value = 38333
$ws.Range("N85").Value = -40829
# Row 102
$ws.Range("H102").Value = 1971.8636
$ws.Range("I102").Value = 1751.6316
$ws.Range("K102").Value = 1751.6316
$ws.Range("M102").Value = -129.6315999999999
# Row 116
$ws.Range("H116").Value = 2326339
$ws.Range("I116").Value = 2326339
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2326339
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -2324045
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 2418.5557
$ws.Range("I122").Value = 2375.5
$ws.Range("J122").Value = 2504.6667
$ws.Range("K122").Value = 7126.5
$ws.Range("L122").Value = 7514.000100000001
$ws.Range("M122").Value = -4676.5
$ws.Range("N122").Value = -12414.0001
# Row 136
$ws.Range("H136").Value = 5937
$ws.Range("I136").Value = 3200.5557
$ws.Range("J136").Value = 9455.286
$ws.Range("K136").Value = 9601.667099999999
$ws.Range("L136").Value = 28365.858
$ws.Range("M136").Value = -7051.667099999999
$ws.Range("N136").Value = -33465.858

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2326339
$ws.Range("I3").Value = 2326339
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2326339
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2326225
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 5799
$ws.Range("J15").Value = 5799
$ws.Range("L15").Value = 5799
$ws.Range("N15").Value = -6139
# Row 17
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 19
$ws.Range("H19").Value = 675.6
$ws.Range("I19").Value = 675.6
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 675.6
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -505.6
$ws.Range("N19").ClearContents()
# Row 24
$ws.Range("H24").Value = 675.6
$ws.Range("I24").Value = 675.6
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 675.6
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -505.6
$ws.Range("N24").ClearContents()
# Row 41
$ws.Range("H41").Value = 28999.666
$ws.Range("J41").Value = 28999.666
$ws.Range("L41").Value = 28999.666
$ws.Range("N41").Value = -29855.666

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 104
$ws.Range("H104").Value = 3515.4707
$ws.Range("J104").Value = 3717.5334
$ws.Range("L104").Value = 11152.6002
$ws.Range("N104").Value = -16394.6002
# Row 105
$ws.Range("H105").Value = 4666.6665
$ws.Range("J105").Value = 4666.6665
$ws.Range("L105").Value = 13999.9995
$ws.Range("N105").Value = -19241.9995
# Row 106
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 24000
$ws.Range("N106").Value = -25892
# Row 107
$ws.Range("H107").Value = 1430.3125
$ws.Range("J107").Value = 1430.3125
$ws.Range("L107").Value = 4290.9375
$ws.Range("N107").Value = -8130.9375
# Row 129
$ws.Range("H129").Value = 81247
$ws.Range("I129").Value = 689
$ws.Range("J129").Value = 145693.4
$ws.Range("K129").Value = 2067
$ws.Range("L129").Value = 437080.2
$ws.Range("M129").Value = 2933
$ws.Range("N129").Value = -447080.2
# Row 131
$ws.Range("H131").Value = 814.8333
$ws.Range("J131").Value = 829.91766
$ws.Range("L131").Value = 2489.75298
$ws.Range("N131").Value = -12569.75298
# Row 141
$ws.Range("H141").Value = 3079.3333
$ws.Range("I141").Value = 3089.25
$ws.Range("K141").Value = 9267.75
$ws.Range("M141").Value = -4087.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996
# Row 83
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984
# Row 102
$ws.Range("H102").Value = 2532.0557
$ws.Range("I102").Value = 2560.2727
$ws.Range("J102").Value = 2487.7144
$ws.Range("K102").Value = 2560.2727
$ws.Range("L102").Value = 2487.7144
$ws.Range("M102").Value = -938.2727
$ws.Range("N102").Value = -5731.7144
# Row 113
$ws.Range("H113").Value = 1150.5
$ws.Range("I113").Value = 1034.3334
$ws.Range("J113").Value = 1266.6666
$ws.Range("K113").Value = 1034.3334
$ws.Range("L113").Value = 1266.6666
$ws.Range("M113").Value = 1135.6666
$ws.Range("N113").Value = -5606.6666
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 2499.5
$ws.Range("J122").Value = 2499.5
$ws.Range("L122").Value = 7498.5
$ws.Range("N122").Value = -12398.5
# Row 126
$ws.Range("H126").Value = 2177639
$ws.Range("I126").Value = 2926325.8
$ws.Range("J126").Value = 145489.14
$ws.Range("K126").Value = 8778977.399999999
$ws.Range("L126").Value = 436467.42
$ws.Range("M126").Value = -8776507.399999999
$ws.Range("N126").Value = -441407.42

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3016.6667
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 1620
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 1620
$ws.Range("M22").Value = -9705
$ws.Range("N22").Value = -2210
# Row 27
$ws.Range("H27").Value = 3016.6667
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 1620
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 1620
$ws.Range("M27").Value = -9893
$ws.Range("N27").Value = -1834
# Row 46
$ws.Range("H46").Value = 2443.375
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
# Row 55
$ws.Range("H55").Value = 453.0909
$ws.Range("I55").Value = 391.8
$ws.Range("J55").Value = 504.16666
$ws.Range("K55").Value = 391.8
$ws.Range("L55").Value = 504.16666
$ws.Range("M55").Value = -218.8
$ws.Range("N55").Value = -850.16666
# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -11996
# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -39984
# Row 132
$ws.Range("H132").Value = 1843.871
$ws.Range("I132").Value = 1663.25
$ws.Range("J132").Value = 2036.5333
$ws.Range("K132").Value = 4989.75
$ws.Range("L132").Value = 6109.5999
$ws.Range("M132").Value = -2459.75
$ws.Range("N132").Value = -11169.5999

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 610.3333
$ws.Range("I100").Value = 373.5
$ws.Range("K100").Value = 747
$ws.Range("M100").Value = -206
# Row 126
$ws.Range("H126").Value = 7318.3105
$ws.Range("I126").Value = 6749.6665
$ws.Range("J126").Value = 8248.817999999999
$ws.Range("K126").Value = 20248.9995
$ws.Range("L126").Value = 24746.454
$ws.Range("M126").Value = -17778.9995
$ws.Range("N126").Value = -29686.454
